# Auto-generated edit script: updates Excel Market Board price/profit
# figures across multiple FFXIV crafting-job sheets (ALC, ARM, BSM, CRP,
# CUL, GSM, LTW, WVR) per the scheduled-runner refresh.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 86
$ws.Range("H86").Value = 2489.75
$ws.Range("I86").Value = 2471.625
$ws.Range("K86").Value = 2471.625
$ws.Range("M86").Value = -1348.625

# Row 89
$ws.Range("H89").Value = 2489.75
$ws.Range("I89").Value = 2471.625
$ws.Range("K89").Value = 12358.125
$ws.Range("M89").Value = -6742.125

# Row 92
$ws.Range("H92").Value = 845.15
$ws.Range("I92").Value = 906.125
$ws.Range("K92").Value = 906.125
$ws.Range("M92").Value = 341.875

# Row 129
$ws.Range("H129").Value = 1183.4166
$ws.Range("J129").Value = 1246.6818
$ws.Range("L129").Value = 3740.0454
$ws.Range("N129").Value = -13740.0454

# Row 132
$ws.Range("H132").Value = 2269.318
$ws.Range("I132").Value = 2210.7144
$ws.Range("J132").Value = 3500
$ws.Range("K132").Value = 6632.1432
$ws.Range("L132").Value = 10500
$ws.Range("M132").Value = -4102.1432
$ws.Range("N132").Value = -15560

# Row 135
$ws.Range("H135").Value = 1178.5714
$ws.Range("I135").Value = 1202
$ws.Range("J135").Value = 1161
$ws.Range("K135").Value = 10818
$ws.Range("L135").Value = 10449
$ws.Range("M135").Value = -8283
$ws.Range("N135").Value = -15519

# Row 138
$ws.Range("H138").Value = 4001.22
$ws.Range("I138").Value = 2126.5
$ws.Range("J138").Value = 4883.4414
$ws.Range("K138").Value = 6379.5
$ws.Range("L138").Value = 14650.3242
$ws.Range("M138").Value = -1239.5
$ws.Range("N138").Value = -24930.3242

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 12427.672
$ws.Range("I32").Value = 13446.852
$ws.Range("J32").Value = 4565.4287
$ws.Range("K32").Value = 13446.852
$ws.Range("L32").Value = 4565.4287
$ws.Range("M32").Value = -13159.852
$ws.Range("N32").Value = -5139.4287

# Row 61
$ws.Range("H61").Value = 2088.4
$ws.Range("I61").Value = 1359.8948
$ws.Range("K61").Value = 1359.8948
$ws.Range("M61").Value = -1147.8948

# Row 88
$ws.Range("H88").Value = 2619.4285
$ws.Range("I88").Value = 2159
$ws.Range("J88").Value = 3233.3333
$ws.Range("K88").Value = 2159
$ws.Range("L88").Value = 3233.3333
$ws.Range("M88").Value = -1753
$ws.Range("N88").Value = -4045.3333

# Row 91
$ws.Range("H91").Value = 2619.4285
$ws.Range("I91").Value = 2159
$ws.Range("J91").Value = 3233.3333
$ws.Range("K91").Value = 2159
$ws.Range("L91").Value = 3233.3333
$ws.Range("M91").Value = -755
$ws.Range("N91").Value = -6041.3333

# Row 132
$ws.Range("H132").Value = 2297.551
$ws.Range("I132").Value = 1758.1538
$ws.Range("J132").Value = 4401.2
$ws.Range("K132").Value = 5274.4614
$ws.Range("L132").Value = 13203.6
$ws.Range("M132").Value = -2744.4614
$ws.Range("N132").Value = -18263.6

# Row 133
$ws.Range("H133").Value = 178812.75
$ws.Range("J133").Value = 178812.75
$ws.Range("L133").Value = 178812.75
$ws.Range("N133").Value = -183872.75

# Row 136
$ws.Range("H136").Value = 2088.4
$ws.Range("I136").Value = 1359.8948
$ws.Range("K136").Value = 4079.6844
$ws.Range("M136").Value = -1529.6844

$ws = $wb.Worksheets.Item("BSM")
# Row 86
$ws.Range("H86").Value = 93966.17999999999
$ws.Range("I86").Value = 3472.1428
$ws.Range("J86").Value = 252330.75
$ws.Range("K86").Value = 3472.1428
$ws.Range("L86").Value = 252330.75
$ws.Range("M86").Value = -2349.1428
$ws.Range("N86").Value = -254576.75

# Row 89
$ws.Range("H89").Value = 93966.17999999999
$ws.Range("I89").Value = 3472.1428
$ws.Range("J89").Value = 252330.75
$ws.Range("K89").Value = 17360.714
$ws.Range("L89").Value = 1261653.75
$ws.Range("M89").Value = -11744.714
$ws.Range("N89").Value = -1272885.75

# Row 134
$ws.Range("H134").Value = 2279.6155
$ws.Range("I134").Value = 1942.7906
$ws.Range("K134").Value = 5828.3718
$ws.Range("M134").Value = -3293.3718

$ws = $wb.Worksheets.Item("CRP")
# Row 58
$ws.Range("H58").Value = 2110.32
$ws.Range("I58").Value = 1832.2667
$ws.Range("J58").Value = 2527.4
$ws.Range("K58").Value = 1832.2667
$ws.Range("L58").Value = 2527.4
$ws.Range("M58").Value = -1629.2667
$ws.Range("N58").Value = -2933.4

# Row 62
$ws.Range("H62").Value = 44325.418
$ws.Range("I62").Value = 57945
$ws.Range("J62").Value = 3466.6667
$ws.Range("K62").Value = 57945
$ws.Range("L62").Value = 3466.6667
$ws.Range("M62").Value = -57321
$ws.Range("N62").Value = -4714.6667

# Row 65
$ws.Range("H65").Value = 44325.418
$ws.Range("I65").Value = 57945
$ws.Range("J65").Value = 3466.6667
$ws.Range("K65").Value = 289725
$ws.Range("L65").Value = 17333.3335
$ws.Range("M65").Value = -286605
$ws.Range("N65").Value = -23573.3335

# Row 132
$ws.Range("H132").Value = 1634.4584
$ws.Range("I132").Value = 1304.6316
$ws.Range("J132").Value = 2887.8
$ws.Range("K132").Value = 3913.8948
$ws.Range("L132").Value = 8663.400000000001
$ws.Range("M132").Value = -1383.8948
$ws.Range("N132").Value = -13723.4

# Row 134
$ws.Range("H134").Value = 1568.5
$ws.Range("I134").Value = 1410.0968
$ws.Range("K134").Value = 4230.2904
$ws.Range("M134").Value = -1695.2904

# Row 136
$ws.Range("H136").Value = 2110.32
$ws.Range("I136").Value = 1832.2667
$ws.Range("J136").Value = 2527.4
$ws.Range("K136").Value = 5496.800099999999
$ws.Range("L136").Value = 7582.200000000001
$ws.Range("M136").Value = -2946.800099999999
$ws.Range("N136").Value = -12682.2

$ws = $wb.Worksheets.Item("CUL")
# Row 6
$ws.Range("H6").Value = 1783.6111
$ws.Range("I6").Value = 650.5
$ws.Range("J6").Value = 1925.25
$ws.Range("K6").Value = 1951.5
$ws.Range("L6").Value = 5775.75
$ws.Range("M6").Value = -1838.5
$ws.Range("N6").Value = -6001.75

# Row 117
$ws.Range("H117").Value = 56372.223
$ws.Range("J117").Value = 72192.86
$ws.Range("L117").Value = 216578.58
$ws.Range("N117").Value = -223462.58

# Row 118
$ws.Range("H118").Value = 2421.75
$ws.Range("I118").Value = 1315
$ws.Range("J118").Value = 2790.6667
$ws.Range("K118").Value = 3945
$ws.Range("L118").Value = 8372.000100000001
$ws.Range("M118").Value = -2702
$ws.Range("N118").Value = -10858.0001

# Row 125
$ws.Range("H125").Value = 3711.111
$ws.Range("I125").Value = 0
$ws.Range("J125").Value = 3711.111
$ws.Range("K125").Value = 0
$ws.Range("L125").Value = 11133.333
$ws.Range("M125").ClearContents()
$ws.Range("N125").Value = -20973.333

# Row 132
$ws.Range("H132").Value = 1495
$ws.Range("I132").Value = 1166.6666
$ws.Range("K132").Value = 10499.9994
$ws.Range("M132").Value = -7969.999400000001

$ws = $wb.Worksheets.Item("GSM")
# Row 70
$ws.Range("H70").Value = 8411.315000000001
$ws.Range("I70").Value = 10661.5
$ws.Range("K70").Value = 10661.5
$ws.Range("M70").Value = -10391.5

# Row 73
$ws.Range("H73").Value = 8411.315000000001
$ws.Range("I73").Value = 10661.5
$ws.Range("K73").Value = 10661.5
$ws.Range("M73").Value = -9725.5

$ws = $wb.Worksheets.Item("LTW")
# Row 122
$ws.Range("H122").Value = 28575266
$ws.Range("I122").Value = 4480
$ws.Range("J122").Value = 40003580
$ws.Range("K122").Value = 13440
$ws.Range("L122").Value = 120010740
$ws.Range("N122").Value = -120015640
$ws.Range("M122").Value = -10990

# Row 136
$ws.Range("H136").Value = 2652.348
$ws.Range("I136").Value = 1633.6
$ws.Range("J136").Value = 4562.5
$ws.Range("K136").Value = 4900.799999999999
$ws.Range("L136").Value = 13687.5
$ws.Range("M136").Value = -2350.799999999999
$ws.Range("N136").Value = -18787.5

$ws = $wb.Worksheets.Item("WVR")
# Row 107
$ws.Range("H107").Value = 773.5
$ws.Range("I107").Value = 774.1111
$ws.Range("J107").Value = 771.6667
$ws.Range("K107").Value = 2322.3333
$ws.Range("L107").Value = 2315.0001
$ws.Range("M107").Value = -402.3332999999998
$ws.Range("N107").Value = -6155.0001
